$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.716.40"
$ws.Range("E2").Value = "  +4.52%  "

$ws.Range("D3").Value = "2.272.53"
$ws.Range("E3").Value = "  +1.73%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.40"
$ws.Range("E5").Value = "  -0.57%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.624"
$ws.Range("E6").Value = "  +0.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.81"
$ws.Range("E7").Value = "  +6.14%  "

$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.426"
$ws.Range("E9").Value = "  +5.52%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.101"
$ws.Range("E10").Value = "  +11.86%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.42"
$ws.Range("E11").Value = "  -0.90%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.83"
$ws.Range("E12").Value = "  +14.29%  "

$ws.Range("E13").Value = "  -0.33%  "

$ws.Range("D14").Value = "2.607.08"
$ws.Range("E14").Value = "  +1.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.62"
$ws.Range("E15").Value = "  +0.76%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.86"
$ws.Range("E16").Value = "  +3.69%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.815"
$ws.Range("E17").Value = "  +1.76%  "

$ws.Range("D18").Value = "2.268.43"
$ws.Range("E18").Value = "  +1.49%  "

$ws.Range("D19").Value = "43.607.71"
$ws.Range("E19").Value = "  +4.42%  "

$ws.Range("D20").Value = "0.0₃0973"
$ws.Range("E20").Value = "  +6.82%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.04"
$ws.Range("E21").Value = "  +0.81%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.08"
$ws.Range("E22").Value = "  -1.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.27"
$ws.Range("E23").Value = "  +0.19%  "

$ws.Range("E24").Value = "  +0.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.47"
$ws.Range("E25").Value = "  +4.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.32"
$ws.Range("E26").Value = "  +0.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.79"
$ws.Range("E27").Value = "  +0.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "171.70"
$ws.Range("E28").Value = "  +1.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.137"
$ws.Range("E29").Value = "  -2.45%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.48"
$ws.Range("E30").Value = "  +2.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.45"
$ws.Range("E31").Value = "  +2.39%  "

$ws.Range("E32").Value = "  +7.77%  "

$ws.Range("E33").Value = "  -0.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0688"
$ws.Range("E34").Value = "  +4.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.11"
$ws.Range("E35").Value = "  +1.98%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.69"
$ws.Range("E36").Value = "  -0.27%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.79"
$ws.Range("E37").Value = "  +5.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.71"
$ws.Range("E38").Value = "  +2.24%  "

$ws.Range("E39").Value = "  -3.47%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0247"
$ws.Range("E40").Value = "  +2.86%  "

$ws.Range("E41").Value = "  -0.29%  "

$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.63"
$ws.Range("E42").Value = "  +20.74%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.34"
$ws.Range("E43").Value = "  -3.81%  "

$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.16"
$ws.Range("E44").Value = "  +2.95%  "

$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0960"
$ws.Range("E45").Value = "  -0.65%  "

$ws.Range("E46").Value = "  -1.37%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "97.02"
$ws.Range("E47").Value = "  -2.21%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.471.94"
$ws.Range("E48").Value = "  -0.25%  "

$ws.Range("B49").Value = "TerraClassic"
$ws.Range("C49").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000211"
$ws.Range("E49").Value = "  -13.30%  "

$ws.Range("B50").Value = "FTXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.37"
$ws.Range("E50").Value = "  -1.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.33"
$ws.Range("E51").Value = "  +3.42%  "
